$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived data for the three remaining rows (sender cluster = MuSCs).
# The old "ECs" sender rows (2-4) are replaced by refreshed "MuSCs" sender
# values, and the trailing duplicate "MuSCs" sender rows (5-7) are removed.
$data = @(
    @("MuSCs","Slitrk2","Ptprs","ECs",   1, 0.3333333333333333, 0.06454900000000001, 0.193647, 1, 1, 3, 1, 1.660421,          4.981262999999999, 0.03714789785507311, 0.03714789785507311, 0.107178515129,    0.964606636161,     0.03714789785507311, 0.03714789785507311),
    @("MuSCs","Slitrk2","Ptprs","FAPs",  1, 0.3333333333333333, 0.06454900000000001, 0.193647, 1, 1, 3, 1, 25.17096033333333, 75.51288099999999, 0.5631392661118858,  0.5631392661118859,  1.624760318556333, 14.622842867007,    0.5631392661118858,  0.5631392661118859),
    @("MuSCs","Slitrk2","Ptprs","MuSCs", 1, 0.3333333333333333, 0.06454900000000001, 0.193647, 1, 1, 3, 1, 17.866195,         53.598585,          0.399712836033041,   0.399712836033041,   1.153245021055,    10.379205189495,    0.399712836033041,   0.399712836033041)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Remove the now-obsolete rows 5:7 (formerly duplicate MuSCs-sender rows).
$ws.Rows("5:7").Delete()
